$wb = $excel.ActiveWorkbook

# Update the CSC curve shape parameter on the "About" sheet (B13),
# following the updated electricity calibration runs.
$ws = $wb.Worksheets.Item("About")
$ws.Range("B13").Value = 0.275
